# Insert a new data row at row 389 (shifts existing rows 389:472 down to 390:473)
# and populate it with the new "Ají" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("389:389").Insert()

$ws.Range("A389").Value = 2
$ws.Range("B389").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C389").Value = "Coquimbo"
$ws.Range("D389").Value = 45211
$ws.Range("E389").Value = 4
$ws.Range("F389").Value = 100112021
$ws.Range("G389").Value = "Ají"
$ws.Range("H389").Value = "Americana (o)"
$ws.Range("I389").Value = "Primera"
$ws.Range("J389").Value = 120
$ws.Range("K389").Value = 35000
$ws.Range("L389").Value = 40000
$ws.Range("M389").Value = 37500
$ws.Range("N389").Value = "$/caja 25 kilos"
$ws.Range("O389").Value = "Provincia de Limarí"
$ws.Range("P389").Value = 1500
$ws.Range("Q389").Value = 25
$ws.Range("R389").Value = "Hortaliza"
